# Fase2 - aumentati semi di zucca per omega-6

$wb = $excel.ActiveWorkbook

# --- Fase2 sheet: increase pumpkin/squash seeds quantity and its %NT ---
$ws2 = $wb.Worksheets.Item("Fase2")
$ws2.Activate()
$ws2.Range("B12").Value = 30
$ws2.Range("E12").Value = 3

# --- Timeline sheet: Fase2 now lasts 13 days instead of 12 ---
$ws6 = $wb.Worksheets.Item("Timeline")
$ws6.Activate()
$ws6.Range("G3").Formula = "=F3+D3/(13/7)"
$ws6.Range("G4").Select()

# Restore Fase2 as the active sheet/selection, matching the saved view state
$ws2.Activate()
$ws2.Range("E13").Select()
